$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRF = 6.84074074074074

for ($r = 28; $r -le 63; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
